$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift rows 11-15 down to 12-16 (copy value+style) to make room for a new
# "Contact" sub-row (the old duplicated Contact row becomes the new
# Jurisdiction row, and everything below it moves down by one).
for ($r = 15; $r -ge 11; $r--) {
    $dest = $r + 1
    $ws.Range("A$r").Copy($ws.Range("A$dest"))
    $ws.Range("B$r").Copy($ws.Range("B$dest"))
}

# Update the metadata values that changed.
$ws.Range("B3").Value = "0.1.7"
$ws.Range("B6").Value = "draft"
$ws.Range("B8").Value = "2024-11-22T12:33:30-06:00"
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

$ws.Range("A11").Value = "Contact"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""

$ws.Range("A13").Value = "Description"
$ws.Range("B13").Value = "SNOMED: Disorders of hematopoietic structure"

$ws.Range("A14").Value = "Purpose"
$ws.Range("B14").Value = ""

$ws.Range("A15").Value = "Copyright"
$ws.Range("B15").Value = ""

$ws.Range("A16").Value = "Immutable"
$ws.Range("B16").Value = "BooleanType[null]"
